$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Title slide: merge "Prof" + ". " runs into a single "Prof. " run.
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$profTr = $s1.Shapes.Item(2).TextFrame.TextRange
$profTr.Characters(1, 6).Text = "Prof. "

# ---------------------------------------------------------------------------
# 2) Duplicate the last slide ("Homework") so a copy with the original
#    content ends up at the very end of the deck, then turn the original
#    (still in place, second-to-last) into the new "Lossy Compression"
#    slide. That reproduces the diff: a brand-new slide about lossy
#    compression is inserted right before the (unchanged) Homework slide.
# ---------------------------------------------------------------------------
$lastIndex = $p.Slides.Count
$homework = $p.Slides.Item($lastIndex)
$homework.Duplicate() | Out-Null

$lossy = $p.Slides.Item($lastIndex)

function EmuToPt($emu) {
    return ($emu / 12700.0) + 0.00002
}

# --- Title: "Lossy" + " Compression" -----------------------------------
$titleTr = $lossy.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Lossy"
$titleTr.InsertAfter(" Compression") | Out-Null

# --- Body placeholder: reposition/resize, then fill with the new text ---
$bodyShape = $lossy.Shapes.Item(2)
$bodyShape.Left = EmuToPt(210553)
$bodyShape.Top = EmuToPt(2603499)
$bodyShape.Width = EmuToPt(12603079)
$bodyShape.Height = EmuToPt(6961605)

$bodyTr = $bodyShape.TextFrame.TextRange

$para1 = "So far, discussed Lossless Compression"
$para2 = "from compressed data, always able to recover the original in full"
$para3 = "To compress even more, could use Lossy Compression"
$para4 = "lose some information when compress, so cannot recover the original"
$para5 = "useful when a decrease in quality is acceptable"
$para6 = "eg: images like JPEG, where quality is degraded to get smaller file size"
$para7 = "eg: music formats like MP3, where removing some sound components that anyway would not be hearable by humans"

$full = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4 + "`r" + $para5 + "`r" + $para6 + "`r" + $para7
$bodyTr.Text = $full

# Paragraph (0-based "lvl") indentation -> COM IndentLevel is 1-based.
$offset = 1
$bodyTr.Characters($offset, $para1.Length).IndentLevel = 1
$offset += $para1.Length + 1
$bodyTr.Characters($offset, $para2.Length).IndentLevel = 2
$offset += $para2.Length + 1
$bodyTr.Characters($offset, $para3.Length).IndentLevel = 1
$offset += $para3.Length + 1
$bodyTr.Characters($offset, $para4.Length).IndentLevel = 2
$offset += $para4.Length + 1
$bodyTr.Characters($offset, $para5.Length).IndentLevel = 2
$offset += $para5.Length + 1
$bodyTr.Characters($offset, $para6.Length).IndentLevel = 2
$offset += $para6.Length + 1
$bodyTr.Characters($offset, $para7.Length).IndentLevel = 2

# Italics for specific sub-strings.
$bodyTr.Characters($full.IndexOf("Lossless Compression") + 1, "Lossless Compression".Length).Font.Italic = $true

$p3Start = $full.IndexOf($para3) + 1
$bodyTr.Characters($p3Start + $para3.IndexOf("Lossy Compression"), "Lossy".Length).Font.Italic = $true
$bodyTr.Characters($p3Start + $para3.IndexOf("Lossy Compression") + "Lossy".Length, " Compression".Length).Font.Italic = $true

$bodyTr.Characters($full.IndexOf($para6) + 1 + $para6.IndexOf("JPEG"), "JPEG".Length).Font.Italic = $true
$bodyTr.Characters($full.IndexOf($para7) + 1 + $para7.IndexOf("MP3"), "MP3".Length).Font.Italic = $true
